$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.9899211372071477
$ws.Cells.Item(2, 3).Value = 0.1900831602652886
$ws.Cells.Item(2, 4).Value = 0.314089528842544
$ws.Cells.Item(2, 6).Value = 1.266215020791975
$ws.Cells.Item(2, 7).Value = 0.002425380193969275
$ws.Cells.Item(2, 10).Value = 0.3480526526176106
$ws.Cells.Item(2, 14).Value = 1.027120027740459
$ws.Cells.Item(2, 15).Value = 2.684521106378014
$ws.Cells.Item(3, 2).Value = 0.8894127075645315
$ws.Cells.Item(3, 3).Value = 0.1659512469055073
$ws.Cells.Item(3, 4).Value = 0.3058457757676223
$ws.Cells.Item(3, 6).Value = 1.257646374425207
$ws.Cells.Item(3, 7).Value = 0.002428501411045452
$ws.Cells.Item(3, 10).Value = 0.3365539059641662
$ws.Cells.Item(3, 14).Value = 1.034545729596829
$ws.Cells.Item(3, 15).Value = 2.677874033397018
$ws.Cells.Item(4, 2).Value = 0.8277868553781218
$ws.Cells.Item(4, 3).Value = 0.1510883087238994
$ws.Cells.Item(4, 4).Value = 0.3009149151776853
$ws.Cells.Item(4, 6).Value = 1.253240974212474
$ws.Cells.Item(4, 7).Value = 0.002430520183762171
$ws.Cells.Item(4, 10).Value = 0.3297037055855867
$ws.Cells.Item(4, 14).Value = 1.039560466073389
$ws.Cells.Item(4, 15).Value = 2.675772641069614
$ws.Cells.Item(5, 2).Value = 0.8026968103210379
$ws.Cells.Item(5, 3).Value = 0.1450203189224339
$ws.Cells.Item(5, 4).Value = 0.2989385378749176
$ws.Cells.Item(5, 6).Value = 1.251660590893081
$ws.Cells.Item(5, 7).Value = 0.002431368661690692
$ws.Cells.Item(5, 10).Value = 0.3269649046592065
$ws.Cells.Item(5, 14).Value = 1.041718604466091
$ws.Cells.Item(5, 15).Value = 2.67541310227486
$ws.Cells.Item(6, 2).Value = 0.7985320509225744
$ws.Cells.Item(6, 3).Value = 0.1440120651267875
$ws.Cells.Item(6, 4).Value = 0.2986123569680359
$ws.Cells.Item(6, 6).Value = 1.251411136126976
$ws.Cells.Item(6, 7).Value = 0.002431511112150106
$ws.Cells.Item(6, 10).Value = 0.3265133091806689
$ws.Cells.Item(6, 14).Value = 1.042083885783804
$ws.Cells.Item(6, 15).Value = 2.67538337496498
$ws.Cells.Item(7, 2).Value = 0.8274483876023169
$ws.Cells.Item(7, 3).Value = 0.151006518720294
$ws.Cells.Item(7, 4).Value = 0.3008881273752309
$ws.Cells.Item(7, 6).Value = 1.253218791132909
$ws.Cells.Item(7, 7).Value = 0.002430531521796231
$ws.Cells.Item(7, 10).Value = 0.3296665558810616
$ws.Cells.Item(7, 14).Value = 1.039589107371384
$ws.Cells.Item(7, 15).Value = 2.67576578213135
$ws.Cells.Item(8, 2).Value = 0.9552485054746853
$ws.Cells.Item(8, 3).Value = 0.1817721681154012
$ws.Cells.Item(8, 4).Value = 0.3112199772188546
$ws.Cells.Item(8, 6).Value = 1.263082721881091
$ws.Cells.Item(8, 7).Value = 0.002426435198400644
$ws.Cells.Item(8, 10).Value = 0.3440442102001811
$ws.Cells.Item(8, 14).Value = 1.029585963221628
$ws.Cells.Item(8, 15).Value = 2.681817601116819
$ws.Cells.Item(9, 2).Value = 1.206513158375628
$ws.Cells.Item(9, 3).Value = 0.241729171712052
$ws.Cells.Item(9, 4).Value = 0.332516306200489
$ws.Cells.Item(9, 6).Value = 1.289234470364107
$ws.Cells.Item(9, 7).Value = 0.002419210660368432
$ws.Cells.Item(9, 10).Value = 0.3739126520952709
$ws.Cells.Item(9, 14).Value = 1.013578854709245
$ws.Cells.Item(9, 15).Value = 2.70944864490852
$ws.Cells.Item(10, 2).Value = 1.391479133488531
$ws.Cells.Item(10, 3).Value = 0.2855414243141752
$ws.Cells.Item(10, 4).Value = 0.3487926481312229
$ws.Cells.Item(10, 6).Value = 1.312628590024502
$ws.Cells.Item(10, 7).Value = 0.002414390496834008
$ws.Cells.Item(10, 10).Value = 0.3968899363916449
$ws.Cells.Item(10, 14).Value = 1.004014243407198
$ws.Cells.Item(10, 15).Value = 2.739440311970327
$ws.Cells.Item(11, 2).Value = 1.475697572403419
$ws.Cells.Item(11, 3).Value = 0.3054192072157775
$ws.Cells.Item(11, 4).Value = 0.356333797319337
$ws.Cells.Item(11, 6).Value = 1.324185611464785
$ws.Cells.Item(11, 7).Value = 0.002412302498366371
$ws.Cells.Item(11, 10).Value = 0.4075700496539127
$ws.Cells.Item(11, 14).Value = 1.000139095005466
$ws.Cells.Item(11, 15).Value = 2.755206640767454
$ws.Cells.Item(12, 2).Value = 1.50759898494033
$ws.Cells.Item(12, 3).Value = 0.3129385890691196
$ws.Cells.Item(12, 4).Value = 0.3592090720117369
$ws.Cells.Item(12, 6).Value = 1.328693978251664
$ws.Cells.Item(12, 7).Value = 0.002411526804135547
$ws.Cells.Item(12, 10).Value = 0.4116472368073403
$ws.Cells.Item(12, 14).Value = 0.9987400498805528
$ws.Cells.Item(12, 15).Value = 2.761483546210314
$ws.Cells.Item(13, 2).Value = 1.500728033111386
$ws.Cells.Item(13, 3).Value = 0.3113195124575725
$ws.Cells.Item(13, 4).Value = 0.3585889596363359
$ws.Cells.Item(13, 6).Value = 1.327717144705261
$ws.Cells.Item(13, 7).Value = 0.002411693198520847
$ws.Cells.Item(13, 10).Value = 0.4107676781172529
$ws.Cells.Item(13, 14).Value = 0.9990383181462192
$ws.Cells.Item(13, 15).Value = 2.760118048844788
$ws.Cells.Item(14, 2).Value = 1.478321933378083
$ws.Cells.Item(14, 3).Value = 0.3060379915890223
$ws.Cells.Item(14, 4).Value = 0.3565699554922048
$ws.Cells.Item(14, 6).Value = 1.324553870217258
$ws.Cells.Item(14, 7).Value = 0.002412238381757303
$ws.Cells.Item(14, 10).Value = 0.4079048230010329
$ws.Cells.Item(14, 14).Value = 1.000022624278159
$ws.Cells.Item(14, 15).Value = 2.755716894842266
$ws.Cells.Item(15, 2).Value = 1.464598777902495
$ws.Cells.Item(15, 3).Value = 0.3028018684239839
$ws.Cells.Item(15, 4).Value = 0.3553358076812572
$ws.Cells.Item(15, 6).Value = 1.322633472282675
$ws.Cells.Item(15, 7).Value = 0.002412574271531106
$ws.Cells.Item(15, 10).Value = 0.4061555247056248
$ws.Cells.Item(15, 14).Value = 1.000634445572032
$ws.Cells.Item(15, 15).Value = 2.753061018742244
$ws.Cells.Item(16, 2).Value = 1.385976689908489
$ws.Cells.Item(16, 3).Value = 0.2842412767238045
$ws.Cells.Item(16, 4).Value = 0.34830256461575
$ws.Cells.Item(16, 6).Value = 1.311891756966858
$ws.Cells.Item(16, 7).Value = 0.0024145290534265
$ws.Cells.Item(16, 10).Value = 0.3961965572361095
$ws.Cells.Item(16, 14).Value = 1.004277065399769
$ws.Cells.Item(16, 15).Value = 2.73845277239576
$ws.Cells.Item(17, 2).Value = 1.337763262420083
$ws.Cells.Item(17, 3).Value = 0.2728412361254868
$ws.Cells.Item(17, 4).Value = 0.3440229126333918
$ws.Cells.Item(17, 6).Value = 1.305536686128079
$ws.Cells.Item(17, 7).Value = 0.002415755017593175
$ws.Cells.Item(17, 10).Value = 0.3901454452598898
$ws.Cells.Item(17, 14).Value = 1.006633541530242
$ws.Cells.Item(17, 15).Value = 2.730035724438864
$ws.Cells.Item(18, 2).Value = 1.310039449858493
$ws.Cells.Item(18, 3).Value = 0.2662793003783577
$ws.Cells.Item(18, 4).Value = 0.3415742666136907
$ws.Cells.Item(18, 6).Value = 1.301967509418589
$ws.Cells.Item(18, 7).Value = 0.002416470020694484
$ws.Cells.Item(18, 10).Value = 0.3866864238271717
$ws.Cells.Item(18, 14).Value = 1.008033714447784
$ws.Cells.Item(18, 15).Value = 2.725394209544106
$ws.Cells.Item(19, 2).Value = 1.300653932051148
$ws.Cells.Item(19, 3).Value = 0.2640567045086755
$ws.Cells.Item(19, 4).Value = 0.3407474148824861
$ws.Cells.Item(19, 6).Value = 1.30077382144114
$ws.Cells.Item(19, 7).Value = 0.002416713804213361
$ws.Cells.Item(19, 10).Value = 0.3855189333283562
$ws.Cells.Item(19, 14).Value = 1.008515482894367
$ws.Cells.Item(19, 15).Value = 2.723856941052219
$ws.Cells.Item(20, 2).Value = 1.342894920373851
$ws.Cells.Item(20, 3).Value = 0.2740553029836406
$ws.Cells.Item(20, 4).Value = 0.3444771547930543
$ws.Cells.Item(20, 6).Value = 1.306204281109231
$ws.Cells.Item(20, 7).Value = 0.002415623491713505
$ws.Cells.Item(20, 10).Value = 0.39078737884887
$ws.Cells.Item(20, 14).Value = 1.006378055285047
$ws.Cells.Item(20, 15).Value = 2.730911051226855
$ws.Cells.Item(21, 2).Value = 1.484902894454081
$ws.Cells.Item(21, 3).Value = 0.3075895194428995
$ws.Cells.Item(21, 4).Value = 0.3571624545832606
$ws.Cells.Item(21, 6).Value = 1.325479415619284
$ws.Cells.Item(21, 7).Value = 0.002412077842060072
$ws.Cells.Item(21, 10).Value = 0.4087448201187271
$ws.Cells.Item(21, 14).Value = 0.9997316538922121
$ws.Cells.Item(21, 15).Value = 2.757001291037909
$ws.Cells.Item(22, 2).Value = 1.57776936533935
$ws.Cells.Item(22, 3).Value = 0.3294598487232179
$ws.Cells.Item(22, 4).Value = 0.3655672640115597
$ws.Cells.Item(22, 6).Value = 1.33884626424387
$ws.Cells.Item(22, 7).Value = 0.002409847868112822
$ws.Cells.Item(22, 10).Value = 0.4206726389245574
$ws.Cells.Item(22, 14).Value = 0.9957864602204722
$ws.Cells.Item(22, 15).Value = 2.775840100163691
$ws.Cells.Item(23, 2).Value = 1.528200075303346
$ws.Cells.Item(23, 3).Value = 0.3177915817272776
$ws.Cells.Item(23, 4).Value = 0.3610710361730582
$ws.Cells.Item(23, 6).Value = 1.331641586726533
$ws.Cells.Item(23, 7).Value = 0.002411030081983337
$ws.Cells.Item(23, 10).Value = 0.4142889646893053
$ws.Cells.Item(23, 14).Value = 0.9978556214225023
$ws.Cells.Item(23, 15).Value = 2.765621515968604
$ws.Cells.Item(24, 2).Value = 1.340574914132901
$ws.Cells.Item(24, 3).Value = 0.2735064479437028
$ws.Cells.Item(24, 4).Value = 0.3442717552046872
$ws.Cells.Item(24, 6).Value = 1.305902198424704
$ws.Cells.Item(24, 7).Value = 0.002415682922627063
$ws.Cells.Item(24, 10).Value = 0.3904970988770486
$ws.Cells.Item(24, 14).Value = 1.006493419179748
$ws.Cells.Item(24, 15).Value = 2.730514700684807
$ws.Cells.Item(25, 2).Value = 1.138473451142261
$ws.Cells.Item(25, 3).Value = 0.2255503241127599
$ws.Cells.Item(25, 4).Value = 0.3266443632966798
$ws.Cells.Item(25, 6).Value = 1.281427754380971
$ws.Cells.Item(25, 7).Value = 0.002421079078218608
$ws.Cells.Item(25, 10).Value = 0.3656519558440436
$ws.Cells.Item(25, 14).Value = 1.017523304407447
$ws.Cells.Item(25, 15).Value = 2.67541310227486
